$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3971766666666667
$ws.Range("H2").Value = 1.19153
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3208016666666667
$ws.Range("N2").Value = 0.962405
$ws.Range("O2").Value = 0.1214789480184765
$ws.Range("P2").Value = 0.1309586095925594
$ws.Range("Q2").Value = 0.1274149366277778
$ws.Range("R2").Value = 1.14673442965
$ws.Range("S2").Value = 0.1214789480184765
$ws.Range("T2").Value = 0.1309586095925594

# Row 3 (Target cluster: FAPs)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3971766666666667
$ws.Range("H3").Value = 1.19153
$ws.Range("N3").Value = 5.205038999999999
$ws.Range("O3").Value = 0.6570026777865272
$ws.Range("P3").Value = 0.7082721622550237
$ws.Range("Q3").Value = 0.6891066799633332
$ws.Range("R3").Value = 6.201960119669999
$ws.Range("S3").Value = 0.6570026777865272
$ws.Range("T3").Value = 0.7082721622550237

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3971766666666667
$ws.Range("H4").Value = 1.19153
$ws.Range("M4").Value = 0.01150933333333333
$ws.Range("N4").Value = 0.034528
$ws.Range("O4").Value = 0.004358274444939456
$ws.Range("P4").Value = 0.004698374252016452
$ws.Range("Q4").Value = 0.00457123864888889
$ws.Range("R4").Value = 0.04114114784
$ws.Range("S4").Value = 0.004358274444939456
$ws.Range("T4").Value = 0.004698374252016452

# Row 5 (Target cluster: MuSCs)
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3971766666666667
$ws.Range("H5").Value = 1.19153
$ws.Range("M5").Value = 0.5734764999999999
$ws.Range("N5").Value = 1.146953
$ws.Range("O5").Value = 0.2171600997500568
$ws.Range("P5").Value = 0.1560708539004004
$ws.Range("Q5").Value = 0.2277714846816667
$ws.Range("R5").Value = 1.36662890809
$ws.Range("S5").Value = 0.2171600997500568
$ws.Range("T5").Value = 0.1560708539004004
